# cryptos list refresh — Tue Jul 25 20:00:23 UTC 2023 (GitHub Actions)
#
# The source sheet stores Price/Volume(1h) as literal text (values such as
# "29.235.39" or "1.859.31" use '.' as a thousands separator, so they are not
# valid numbers, and "  +0.36%  " is padded text, not a real percentage).
# When Excel's Range.Value setter is handed a string that DOES parse as a
# plain number (e.g. "237.49", "1.000"), it silently converts the cell to a
# Number and can lose formatting (trailing zeros, etc.). To keep those cells
# as plain text — exactly like the rest of the column — we prefix such
# values with a leading apostrophe, Excel's standard "force text" marker;
# Excel strips the marker and stores the clean text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Text
    )
    # If the trimmed text would parse as a plain number, Excel's Value setter
    # auto-converts the cell to a Number (dropping e.g. trailing zeros), so
    # prefix it with an apostrophe — Excel's "force text" marker — to keep it
    # as literal text, matching every other cell in the column.
    $trimmed = $Text.Trim()
    if ($trimmed -match '^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$') {
        $ws.Range($CellRef).Value = "'" + $Text
    } else {
        $ws.Range($CellRef).Value = $Text
    }
}

# --- Rows 41/42: two coins swapped rank/position, each gets new price & volume ---
Set-TextValue "B41" "Maker"
Set-TextValue "C41" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D41" "1.137.27"
Set-TextValue "E41" "  +5.02%  "

Set-TextValue "B42" "TrustWalletToken"
Set-TextValue "C42" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D42" "0.9187"
Set-TextValue "E42" "  -2.34%  "

# --- Price / Volume(1h) refresh for the remaining rows ---
Set-TextValue "D2" "29.235.27"
Set-TextValue "E2" "  +0.40%  "

Set-TextValue "D3" "1.859.12"
Set-TextValue "E3" "  +0.53%  "

Set-TextValue "E4" "  +0.06%  "

Set-TextValue "D5" "0.7010"
Set-TextValue "E5" "  -0.27%  "

Set-TextValue "D6" "237.49"
Set-TextValue "E6" "  -0.14%  "

Set-TextValue "E7" "  +0.05%  "

Set-TextValue "D8" "0.08256"
Set-TextValue "E8" "  +9.87%  "

Set-TextValue "E9" "  +0.00%  "

Set-TextValue "D10" "23.21"
Set-TextValue "E10" "  -0.49%  "

Set-TextValue "D11" "0.08182"
Set-TextValue "E11" "  +0.80%  "

Set-TextValue "D12" "1.859.04"
Set-TextValue "E12" "  +0.70%  "

Set-TextValue "D13" "5.170"
Set-TextValue "E13" "  -0.84%  "

Set-TextValue "D14" "0.7111"
Set-TextValue "E14" "  -1.88%  "

Set-TextValue "D15" "89.04"

Set-TextValue "D16" "29.253.64"
Set-TextValue "E16" "  +0.58%  "

Set-TextValue "D17" "5.772"
Set-TextValue "E17" "  +0.04%  "

Set-TextValue "D18" "0.000007837"
Set-TextValue "E18" "  +2.42%  "

Set-TextValue "D19" "13.34"
Set-TextValue "E19" "  +2.24%  "

Set-TextValue "D20" "236.65"
Set-TextValue "E20" "  -0.63%  "

Set-TextValue "E21" "  +0.06%  "

Set-TextValue "D22" "2.113.09"
Set-TextValue "E22" "  +1.34%  "

Set-TextValue "E23" "  +0.00%  "

Set-TextValue "D24" "7.438"
Set-TextValue "E24" "  -1.30%  "

Set-TextValue "D25" "161.76"
Set-TextValue "E25" "  -0.11%  "

Set-TextValue "D26" "8.974"
Set-TextValue "E26" "  +0.04%  "

Set-TextValue "D27" "0.1441"
Set-TextValue "E27" "  -0.96%  "

Set-TextValue "D28" "18.09"
Set-TextValue "E28" "  +0.46%  "

Set-TextValue "E29" "  +1.26%  "

Set-TextValue "E30" "  +3.47%  "

Set-TextValue "D31" "1.482"
Set-TextValue "E31" "  -0.67%  "

Set-TextValue "D32" "4.392"
Set-TextValue "E32" "  -3.01%  "

Set-TextValue "D33" "4.061"
Set-TextValue "E33" "  +1.89%  "

Set-TextValue "D34" "0.05201"
Set-TextValue "E34" "  +1.17%  "

Set-TextValue "D35" "1.168"
Set-TextValue "E35" "  -1.52%  "

Set-TextValue "D36" "0.7069"

Set-TextValue "D37" "1.003"
Set-TextValue "E37" "  -3.04%  "

Set-TextValue "D38" "2.671"
Set-TextValue "E38" "  +0.39%  "

Set-TextValue "D39" "0.01846"
Set-TextValue "E39" "  -1.46%  "

Set-TextValue "D40" "2.723"
Set-TextValue "E40" "  +1.65%  "

Set-TextValue "D43" "5.927"
Set-TextValue "E43" "  -0.95%  "

Set-TextValue "E44" "  -0.15%  "

Set-TextValue "D45" "70.58"
Set-TextValue "E45" "  +1.05%  "

Set-TextValue "D46" "1.000"
Set-TextValue "E46" "  -0.01%  "

Set-TextValue "E47" "  +0.14%  "

Set-TextValue "D48" "1.769"
Set-TextValue "E48" "  +1.55%  "

Set-TextValue "D49" "2.010.31"
Set-TextValue "E49" "  +1.63%  "

Set-TextValue "D50" "9.158"
Set-TextValue "E50" "  +0.06%  "

Set-TextValue "D51" "6.963"
Set-TextValue "E51" "  -1.06%  "
